# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the active sheet, per the crypto data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.576.73"
$ws.Range("E2").Value = "  -5.89%  "
$ws.Range("D3").Value = "3.261.42"
$ws.Range("E3").Value = "  -6.72%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "554.89"
$ws.Range("E5").Value = "  -4.24%  "
$ws.Range("D6").Value = "'182.90"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -4.05%  "
$ws.Range("D9").Value = "3.254.36"
$ws.Range("E9").Value = "  -6.55%  "
$ws.Range("E10").Value = "  -10.76%  "
$ws.Range("D11").Value = "'0.580"
$ws.Range("E11").Value = "  -6.41%  "
$ws.Range("D12").Value = "46.72"
$ws.Range("E12").Value = "  -9.12%  "
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  -7.74%  "
$ws.Range("D14").Value = "8.58"
$ws.Range("E14").Value = "  -6.20%  "
$ws.Range("D15").Value = "632.15"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "3.788.59"
$ws.Range("E16").Value = "  -6.37%  "
$ws.Range("D17").Value = "17.94"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "65.584.69"
$ws.Range("E18").Value = "  -5.97%  "
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").Value = "3.259.36"
$ws.Range("E20").Value = "  -6.83%  "
$ws.Range("E21").Value = "  -8.78%  "
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("D23").Value = "18.13"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "106.72"
$ws.Range("E24").Value = "  +8.11%  "
$ws.Range("D25").Value = "4.87"
$ws.Range("E25").Value = "  -8.17%  "
$ws.Range("E26").Value = "  -7.99%  "
$ws.Range("E27").Value = "  -7.77%  "
$ws.Range("D28").Value = "9.46"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("D29").Value = "8.61"
$ws.Range("E29").Value = "  -7.87%  "
$ws.Range("D30").Value = "30.06"
$ws.Range("E30").Value = "  -8.08%  "
$ws.Range("D31").Value = "3.91"
$ws.Range("E31").Value = "  -8.74%  "
$ws.Range("E32").Value = "  -7.72%  "
$ws.Range("D33").Value = "10.97"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("D35").Value = "'57.60"
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("D36").Value = "3.730.57"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "521.03"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  -5.95%  "
$ws.Range("D40").Value = "0.0₃0730"
$ws.Range("E40").Value = "  -7.40%  "
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  -7.98%  "
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -5.10%  "
$ws.Range("D44").Value = "32.63"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("D45").Value = "0.333"
$ws.Range("E45").Value = "  -10.56%  "
$ws.Range("E46").Value = "  -7.15%  "
$ws.Range("D47").Value = "3.19"
$ws.Range("E47").Value = "  -5.74%  "
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("E49").Value = "  -9.63%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  +0.81%  "
